$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.027.90'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '1.863.74'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  +0.42%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '312.43'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5091'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3851'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08274'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -8.53%  '
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '41.47'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.45%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.228'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.21%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.58'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('D14').Value = '1.856.33'
$ws.Range('E14').Value = '  -1.05%  '
$ws.Range('E15').Value = '  -0.63%  '
$ws.Range('E16').Value = '  +0.33%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '90.83'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06644'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('E20').Value = '  -2.50%  '
$ws.Range('E21').Value = '  +0.36%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.041'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.56%  '
$ws.Range('D23').Value = '28.041.83'
$ws.Range('E23').Value = '  -0.41%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.09'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.37%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.229'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.45%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.533'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.59%  '
$ws.Range('D27').Value = '2.073.08'
$ws.Range('E27').Value = '  -0.90%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '157.67'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.39%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '20.54'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.44%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '125.55'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.08%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.1060'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.39%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.036'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.60%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.900'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +5.09%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.596'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.425'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.45%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.06553'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02422'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.49%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2176'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.88%  '
$ws.Range('E39').Value = '  -0.82%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.6487'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.21%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.997'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.40%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.221'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -5.52%  '
$ws.Range('E43').Value = '  -2.77%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6128'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.51%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.10'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.49%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.292'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.26%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.652'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.32%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.011'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.24%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.211'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.05%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '120.16'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.01%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '78.52'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.79%  '
